# Updates cryptocurrency price (D) and volume-change (E) figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-detects numeric-looking strings and coerces them to numbers, which would
# drop significant trailing zeros (e.g. "60.00" -> 60) or reformat dotted-thousands
# values. To guarantee every cell keeps the exact text from the source feed, each value
# is written as a quoted text formula, then immediately flattened back to a plain value
# via copy / paste-special, without touching cell styles/number formats.
function Set-TextValue([string]$addr, [string]$val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = "=""" + $escaped + """"
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue "D2" "36.465.11"
Set-TextValue "E2" "  -0.04%  "

Set-TextValue "D3" "1.954.09"
Set-TextValue "E3" "  +0.60%  "

Set-TextValue "E4" "  +0.07%  "

Set-TextValue "D5" "243.45"
Set-TextValue "E5" "  +0.12%  "

Set-TextValue "E6" "  +2.71%  "

Set-TextValue "D7" "60.00"
Set-TextValue "E7" "  +5.09%  "

Set-TextValue "E8" "  -0.02%  "

Set-TextValue "E9" "  +3.93%  "

Set-TextValue "D10" "0.0785"
Set-TextValue "E10" "  -2.74%  "

Set-TextValue "D11" "0.103"
Set-TextValue "E11" "  +0.71%  "

Set-TextValue "D12" "14.10"

Set-TextValue "D13" "0.839"
Set-TextValue "E13" "  +4.22%  "

Set-TextValue "D14" "2.240.42"
Set-TextValue "E14" "  +0.46%  "

Set-TextValue "D15" "21.48"
Set-TextValue "E15" "  -1.07%  "

Set-TextValue "D16" "5.25"
Set-TextValue "E16" "  +1.57%  "

Set-TextValue "D17" "1.950.93"
Set-TextValue "E17" "  +0.43%  "

Set-TextValue "D18" "36.436.47"
Set-TextValue "E18" "  -0.06%  "

Set-TextValue "D19" "69.07"
Set-TextValue "E19" "  -0.27%  "

Set-TextValue "D20" "0.0₃0851"
Set-TextValue "E20" "  -0.28%  "

Set-TextValue "D21" "229.14"
Set-TextValue "E21" "  +0.76%  "

Set-TextValue "D22" "5.06"
Set-TextValue "E22" "  +2.11%  "

Set-TextValue "E23" "  -0.01%  "

Set-TextValue "E24" "  +2.13%  "

Set-TextValue "E25" "  +2.96%  "

Set-TextValue "D26" "0.142"
Set-TextValue "E26" "  +6.05%  "

Set-TextValue "D27" "9.12"
Set-TextValue "E27" "  -0.52%  "

Set-TextValue "D28" "160.44"
Set-TextValue "E28" "  +0.30%  "

Set-TextValue "D29" "19.21"
Set-TextValue "E29" "  +0.22%  "

Set-TextValue "D30" "1.31"
Set-TextValue "E30" "  +20.19%  "

Set-TextValue "E31" "  +2.17%  "

Set-TextValue "D32" "4.77"
Set-TextValue "E32" "  +3.06%  "

Set-TextValue "D33" "0.0610"
Set-TextValue "E33" "  -0.86%  "

Set-TextValue "E34" "  +7.08%  "

Set-TextValue "E35" "  +0.08%  "

Set-TextValue "D36" "2.25"
Set-TextValue "E36" "  +2.93%  "

Set-TextValue "E37" "  +3.49%  "

Set-TextValue "E38" "  -0.19%  "

Set-TextValue "E39" "  -11.39%  "

Set-TextValue "D40" "0.0962"
Set-TextValue "E40" "  -2.40%  "

Set-TextValue "E42" "  +1.57%  "

Set-TextValue "D43" "0.0209"
Set-TextValue "E43" "  +0.55%  "

Set-TextValue "E44" "  +0.46%  "

Set-TextValue "D45" "1.358.97"
Set-TextValue "E45" "  +1.25%  "

Set-TextValue "D46" "88.49"
Set-TextValue "E46" "  +2.76%  "

Set-TextValue "D47" "1.02"
Set-TextValue "E47" "  -0.02%  "

Set-TextValue "D48" "7.18"
Set-TextValue "E48" "  +0.92%  "

Set-TextValue "E49" "  -0.07%  "

Set-TextValue "D50" "46.00"
Set-TextValue "E50" "  +6.68%  "

Set-TextValue "D51" "2.136.11"
Set-TextValue "E51" "  +0.71%  "

$excel.CutCopyMode = $false
